# edit.ps1 - applies the "Game changes + Report progress" commit:
#  - Insert a new "Scene Flow" slide at position 5 (existing slides 5-9 shift down to 6-10)
#  - Append a new "References" slide at the end (position 11)
#  - Update the cached footer date field text 31/03/2025 -> 01/04/2025 on the
#    slide master and all slide layouts

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Insert the new "Scene Flow" slide at position 5
# ---------------------------------------------------------------------------
$titleAndContent = $p.SlideMaster.CustomLayouts.Item(2)

$sceneFlow = $p.Slides.AddSlide(5, $titleAndContent)

$sfTitle = $sceneFlow.Shapes.Item(1).TextFrame.TextRange
$sfTitle.Text = "Scene Flow"
$sfTitle.Font.Bold = $true

$sfBody = $sceneFlow.Shapes.Item(2).TextFrame.TextRange
$sfBody.Text = "MainMenu " + [char]8594 + " Level1 " + [char]8594 + " Loading " + [char]8594 + " Level2 " + [char]8594 + " Loading " + [char]8594 + " Level3 " + [char]8594 + " EndScreen " + [char]8594 + " MainMenu`r"
$sfBody.Font.Bold = $true

$sfPara1 = $sfBody.Paragraphs(1)
$sfPara1.ParagraphFormat.Bullet.Visible = $true
$sfPara1.ParagraphFormat.Bullet.Character = 8226
$sfPara1.ParagraphFormat.Bullet.Font.Name = "Arial"

$sfPara2 = $sfBody.Paragraphs(2)
$sfPara2.Font.Bold = $false

# ---------------------------------------------------------------------------
# 2. Append the new "References" slide at the end
# ---------------------------------------------------------------------------
$refs = $p.Slides.AddSlide($p.Slides.Count + 1, $titleAndContent)

$refTitle = $refs.Shapes.Item(1).TextFrame.TextRange
$refTitle.Text = "References "

$refBody = $refs.Shapes.Item(2).TextFrame.TextRange
$link1 = "https://parents.actionforchildren.org.uk/home-family-life/technology/video-game-age-restrictions/"
$link2 = "https://openr.co/unveiling-the-demographics-of-call-of-dutys-target-audience/"
$refBody.Text = "1. " + $link1 + "`r2. " + $link2 + "`r"

$p1 = $refBody.Paragraphs(1)
$p1.Characters(4, $link1.Length).ActionSettings.Item(1).Hyperlink.Address = $link1

$p2 = $refBody.Paragraphs(2)
$p2.Characters(4, $link2.Length).ActionSettings.Item(1).Hyperlink.Address = $link2

# ---------------------------------------------------------------------------
# 3. Update the cached date-field text on the master and every layout
# ---------------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "31/03/2025") {
                $sh.TextFrame.TextRange.Text = "01/04/2025"
            }
        }
    }
}

Update-DateShape $p.SlideMaster.Shapes

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DateShape $p.SlideMaster.CustomLayouts.Item($li).Shapes
}
